$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered roster for rows 6-14 (Name, Position, Team)
$data = @(
    @("Patrick Williams", "PF", "Chicago Bulls"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Goga Bitadze", "C", "Orlando Magic"),
    @("Brandin Podziemski", "SG", "Golden State Warriors"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers")
)

$r = 6
foreach ($row in $data) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $r++
}
